$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 35, shifting existing rows 35-83 down to 36-84.
$ws.Rows.Item(35).Insert()

# Populate the newly inserted row 35 with the new weekly record.
$ws.Cells.Item(35, 1).Value = 8
$ws.Cells.Item(35, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(35, 3).Value = "Coquimbo"
$ws.Cells.Item(35, 4).Value = 44467
$ws.Cells.Item(35, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(35, 5).Value = 4
$ws.Cells.Item(35, 6).Value = 100112044
$ws.Cells.Item(35, 7).Value = "Perejil"
$ws.Cells.Item(35, 8).Value = "Sin especificar"
$ws.Cells.Item(35, 9).Value = "Primera"
$ws.Cells.Item(35, 10).Value = 3100
$ws.Cells.Item(35, 11).Value = 1500
$ws.Cells.Item(35, 12).Value = 2000
$ws.Cells.Item(35, 13).Value = 1750
$ws.Cells.Item(35, 14).Value = "`$/atado 1 a 1,5 kilos"
$ws.Cells.Item(35, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(35, 16).Value = 1167
$ws.Cells.Item(35, 17).Value = 1.5
$ws.Cells.Item(35, 18).Value = "Hortaliza"
